$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.262.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'1.792.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.59%  '
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").Value = "'325.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.82%  '
$ws.Range("E6").Value = '  +0.32%  '
$ws.Range("D7").Value = "'0.4455"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +13.56%  '
$ws.Range("D8").Value = "'0.3731"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +9.79%  '
$ws.Range("D9").Value = "'44.57"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.69%  '
$ws.Range("E10").Value = '  +1.86%  '
$ws.Range("D11").Value = "'0.07499"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.58%  '
$ws.Range("D12").Value = "'22.56"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.39%  '
$ws.Range("D13").Value = "'1.002"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.22%  '
$ws.Range("D14").Value = "'6.272"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.79%  '
$ws.Range("D15").Value = "'7.523"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.68%  '
$ws.Range("D16").Value = "'1.789.99"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.52%  '
$ws.Range("D17").Value = "'0.00001088"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.38%  '
$ws.Range("D18").Value = "'0.06737"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.85%  '
$ws.Range("D19").Value = "'80.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.54%  '
$ws.Range("E20").Value = '  +0.24%  '
$ws.Range("D21").Value = "'17.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.13%  '
$ws.Range("D22").Value = "'6.311"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.05%  '
$ws.Range("D23").Value = "'28.245.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.64%  '
$ws.Range("D24").Value = "'11.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.27%  '
$ws.Range("D25").Value = "'2.420"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.56%  '
$ws.Range("D26").Value = "'20.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.07%  '
$ws.Range("D27").Value = "'151.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.57%  '
$ws.Range("D28").Value = "'2.354"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.13%  '
$ws.Range("D29").Value = "'1.991.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.74%  '
$ws.Range("D30").Value = "'132.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.14%  '
$ws.Range("D31").Value = "'1.222"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.90%  '
$ws.Range("D32").Value = "'4.023"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.38%  '
$ws.Range("D33").Value = "'5.796"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.81%  '
$ws.Range("D34").Value = "'0.09376"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.24%  '
$ws.Range("D35").Value = "'0.2338"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +10.45%  '
$ws.Range("D36").Value = "'12.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.53%  '
$ws.Range("D37").Value = "'0.06328"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.87%  '
$ws.Range("D38").Value = "'0.02333"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.67%  '
$ws.Range("D39").Value = "'5.149"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.13%  '
$ws.Range("D40").Value = "'0.6533"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.17%  '
$ws.Range("D41").Value = "'8.308"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.85%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = "'1.208"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.14%  '
$ws.Range("B43").Value = 'WEMIXTOKEN'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D43").Value = "'1.467"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.87%  '
$ws.Range("D44").Value = "'1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.26%  '
$ws.Range("D45").Value = "'14.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.56%  '
$ws.Range("D46").Value = "'0.6070"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.82%  '
$ws.Range("D47").Value = "'3.777"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.29%  '
$ws.Range("D48").Value = "'129.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.25%  '
$ws.Range("D49").Value = "'2.021"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.02%  '
$ws.Range("D50").Value = "'0.07115"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.36%  '
$ws.Range("D51").Value = "'1.155"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.58%  '
